# Insert a new price record before row 110 (a new daily reading), which
# pushes the existing rows 110-149 down to 111-150.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(110).Insert()

$ws.Range("A110").Value = 4
$ws.Range("B110").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C110").Value = "Los Lagos"
$ws.Range("D110").Value = 44468
$ws.Range("E110").Value = 10
$ws.Range("F110").Value = 100112043
$ws.Range("G110").Value = "Pepino ensalada"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 150
$ws.Range("K110").Value = 20000
$ws.Range("L110").Value = 20000
$ws.Range("M110").Value = 20000
$ws.Range("N110").Value = "$/caja 60 unidades"
$ws.Range("O110").Value = "Región de Arica y Parinacota"
$ws.Range("P110").Value = 333
$ws.Range("Q110").Value = 60
$ws.Range("R110").Value = "Hortaliza"
